$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("N2-1_m")
$ws.Activate()

# Overwrite formula cells H6 and H21 with a literal value (as a user typing
# over them would), breaking the shared formula group and triggering a
# workbook-wide recalculation.
$ws.Range("H6").Value = -0.3
$ws.Range("H21").Value = -0.3

# Leave the final selection on H22, matching the saved view state.
$ws.Range("H22").Select()

# Break the (now-orphaned) external workbook link.
$links = $wb.LinkSources(1)
if ($links) {
    $wb.BreakLink($links, 1)
}
